$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "31.276.47"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +2.64%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.980.70"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +5.15%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9984"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7958"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +68.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "251.93"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.35%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9962"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.44%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3406"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +18.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "25.66"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +15.63%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06934"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +7.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8408"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +16.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08098"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +4.20%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "101.81"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +6.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.978.26"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +4.98%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.487"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +6.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "274.17"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -2.32%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "31.244.48"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.54%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.93"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +6.84%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007862"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +5.45%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.242.46"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +4.98%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.672"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +7.90%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9952"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.55%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9990"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.826"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +8.69%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1570"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +63.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.631"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +6.38%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "166.15"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.31%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.55"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +3.89%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.232"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +18.31%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.559"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +6.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.349"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.543"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +6.66%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.334"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +4.87%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05187"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +7.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.215"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +8.14%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7413"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +7.17%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.775"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.26%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01980"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +5.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.904"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.92%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.589"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +6.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "78.47"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +5.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4667"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +9.58%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.072"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +5.75%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8539"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +3.35%  "
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "105.54"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +4.42%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9949"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.999"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +3.86%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.476"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +7.54%  "
$ws.Range("B49").Value = "Decentraland"
$ws.Range("C49").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4265"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +8.59%  "
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.35"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +3.40%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "931.33"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +2.98%  "
